$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 46
$ws.Range("H46").Value = 12503088
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 12503088
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 37509264
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -37509502

# Row 47
$ws.Range("H47").Value = 15000
$ws.Range("I47").Value = 0
$ws.Range("K47").Value = 0
$ws.Range("M47").ClearContents()

# Row 48
$ws.Range("H48").Value = 2485.7144
$ws.Range("J48").Value = 2485.7144
$ws.Range("L48").Value = 7457.1432
$ws.Range("N48").Value = -8041.1432

# Row 56
$ws.Range("H56").Value = 2485.7144
$ws.Range("J56").Value = 2485.7144
$ws.Range("L56").Value = 7457.1432
$ws.Range("N56").Value = -8525.143199999999

# Row 60
$ws.Range("H60").Value = 12503088
$ws.Range("I60").Value = 0
$ws.Range("J60").Value = 12503088
$ws.Range("K60").Value = 0
$ws.Range("L60").Value = 37509264
$ws.Range("M60").ClearContents()
$ws.Range("N60").Value = -37510232

# Row 70
$ws.Range("H70").Value = 1859.7
$ws.Range("J70").Value = 2085.0715
$ws.Range("L70").Value = 6255.2145
$ws.Range("N70").Value = -6795.2145

# Row 73
$ws.Range("H73").Value = 1859.7
$ws.Range("J73").Value = 2085.0715
$ws.Range("L73").Value = 6255.2145
$ws.Range("N73").Value = -8127.2145

# Row 134
$ws.Range("H134").Value = 66508.46000000001
$ws.Range("J134").Value = 66508.46000000001
$ws.Range("L134").Value = 66508.46000000001
$ws.Range("N134").Value = -76648.46000000001

# Row 137
$ws.Range("H137").Value = 1152642.8
$ws.Range("I137").Value = 4891.727
$ws.Range("J137").Value = 1854046
$ws.Range("K137").Value = 14675.181
$ws.Range("L137").Value = 5562138
$ws.Range("M137").Value = -12125.181
$ws.Range("N137").Value = -5567238

$ws = $wb.Worksheets.Item("ARM")
# Row 2
$ws.Range("H2").Value = 1002.86664
$ws.Range("I2").Value = 1002.86664
$ws.Range("J2").Value = 0
$ws.Range("K2").Value = 1002.86664
$ws.Range("L2").Value = 0
$ws.Range("M2").Value = -889.86664
$ws.Range("N2").ClearContents()

# Row 32
$ws.Range("H32").Value = 24723.043
$ws.Range("I32").Value = 32901.03
$ws.Range("K32").Value = 32901.03
$ws.Range("M32").Value = -32614.03

# Row 45
$ws.Range("H45").Value = 1682.925
$ws.Range("I45").Value = 1646.5358
$ws.Range("J45").Value = 1767.8334
$ws.Range("K45").Value = 1646.5358
$ws.Range("L45").Value = 1767.8334
$ws.Range("M45").Value = -1269.5358
$ws.Range("N45").Value = -2521.8334

# Row 98
$ws.Range("H98").Value = 34750
$ws.Range("J98").Value = 34750
$ws.Range("L98").Value = 34750
$ws.Range("N98").Value = -40740

# Row 116
$ws.Range("H116").Value = 1002.86664
$ws.Range("I116").Value = 1002.86664
$ws.Range("J116").Value = 0
$ws.Range("K116").Value = 1002.86664
$ws.Range("L116").Value = 0
$ws.Range("M116").Value = 1291.13336
$ws.Range("N116").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
# Row 3
$ws.Range("H3").Value = 1002.86664
$ws.Range("I3").Value = 1002.86664
$ws.Range("J3").Value = 0
$ws.Range("K3").Value = 1002.86664
$ws.Range("L3").Value = 0
$ws.Range("M3").Value = -888.86664
$ws.Range("N3").ClearContents()

# Row 105
$ws.Range("H105").Value = 6843.3335
$ws.Range("I105").Value = 8265
$ws.Range("J105").Value = 4000
$ws.Range("K105").Value = 8265
$ws.Range("L105").Value = 4000
$ws.Range("M105").Value = -6518
$ws.Range("N105").Value = -7494

# Row 107
$ws.Range("H107").Value = 3228.5715
$ws.Range("I107").Value = 3000
$ws.Range("J107").Value = 3800
$ws.Range("K107").Value = 3000
$ws.Range("L107").Value = 3800
$ws.Range("M107").Value = -1080
$ws.Range("N107").Value = -7640

$ws = $wb.Worksheets.Item("CRP")
# Row 62
$ws.Range("H62").Value = 3436.875
$ws.Range("I62").Value = 3499.1667
$ws.Range("J62").Value = 3250
$ws.Range("K62").Value = 3499.1667
$ws.Range("L62").Value = 3250
$ws.Range("M62").Value = -2875.1667
$ws.Range("N62").Value = -4498

# Row 65
$ws.Range("H65").Value = 3436.875
$ws.Range("I65").Value = 3499.1667
$ws.Range("J65").Value = 3250
$ws.Range("K65").Value = 17495.8335
$ws.Range("L65").Value = 16250
$ws.Range("M65").Value = -14375.8335
$ws.Range("N65").Value = -22490

# Row 94
$ws.Range("H94").Value = 804.8
$ws.Range("I94").Value = 674.6667
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 674.6667
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = -223.6667
$ws.Range("N94").Value = -1902

# Row 107
$ws.Range("H107").Value = 948.5625
$ws.Range("I107").Value = 948.5625
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 948.5625
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = 971.4375
$ws.Range("N107").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
# Row 11
$ws.Range("H11").Value = 302
$ws.Range("I11").Value = 199
$ws.Range("J11").Value = 370.66666
$ws.Range("K11").Value = 597
$ws.Range("L11").Value = 1111.99998
$ws.Range("M11").Value = -457
$ws.Range("N11").Value = -1391.99998

# Row 64
$ws.Range("H64").Value = 3937.5
$ws.Range("I64").Value = 750
$ws.Range("J64").Value = 5000
$ws.Range("K64").Value = 2250
$ws.Range("L64").Value = 15000
$ws.Range("M64").Value = -1980
$ws.Range("N64").Value = -15540

# Row 67
$ws.Range("H67").Value = 3937.5
$ws.Range("I67").Value = 750
$ws.Range("J67").Value = 5000
$ws.Range("K67").Value = 2250
$ws.Range("L67").Value = 15000
$ws.Range("M67").Value = -1314
$ws.Range("N67").Value = -16872

# Row 68
$ws.Range("H68").Value = 2644.7463
$ws.Range("I68").Value = 923.37933
$ws.Range("J68").Value = 3958.4211
$ws.Range("K68").Value = 2770.13799
$ws.Range("L68").Value = 11875.2633
$ws.Range("M68").Value = -1959.13799
$ws.Range("N68").Value = -13497.2633

# Row 71
$ws.Range("H71").Value = 2644.7463
$ws.Range("I71").Value = 923.37933
$ws.Range("J71").Value = 3958.4211
$ws.Range("K71").Value = 8310.41397
$ws.Range("L71").Value = 35625.7899
$ws.Range("M71").Value = -4254.41397
$ws.Range("N71").Value = -43737.7899

# Row 98
$ws.Range("H98").Value = 390
$ws.Range("I98").Value = 390
$ws.Range("J98").Value = 0
$ws.Range("K98").Value = 1170
$ws.Range("L98").Value = 0
$ws.Range("M98").Value = 328
$ws.Range("N98").ClearContents()

# Row 102
$ws.Range("H102").Value = 5000
$ws.Range("I102").Value = 0
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 0
$ws.Range("L102").Value = 15000
$ws.Range("M102").ClearContents()
$ws.Range("N102").Value = -19868

# Row 123
$ws.Range("H123").Value = 1210
$ws.Range("I123").Value = 1210
$ws.Range("K123").Value = 3630
$ws.Range("M123").Value = -1180

# Row 129
$ws.Range("H129").Value = 1292.3823
$ws.Range("I129").Value = 647.5
$ws.Range("J129").Value = 1644.1364
$ws.Range("K129").Value = 1942.5
$ws.Range("L129").Value = 4932.4092
$ws.Range("M129").Value = 3057.5
$ws.Range("N129").Value = -14932.4092

# Row 140
$ws.Range("H140").Value = 2414.6897
$ws.Range("I140").Value = 2353.125
$ws.Range("J140").Value = 2490.4614
$ws.Range("K140").Value = 7059.375
$ws.Range("L140").Value = 7471.3842
$ws.Range("M140").Value = -1879.375
$ws.Range("N140").Value = -17831.3842

$ws = $wb.Worksheets.Item("GSM")
# Row 107
$ws.Range("H107").Value = 807.1429000000001
$ws.Range("I107").Value = 212.75
$ws.Range("J107").Value = 1599.6666
$ws.Range("K107").Value = 212.75
$ws.Range("L107").Value = 1599.6666
$ws.Range("M107").Value = 1707.25
$ws.Range("N107").Value = -5439.6666

# Row 141
$ws.Range("H141").Value = 37232.25
$ws.Range("J141").Value = 37232.25
$ws.Range("L141").Value = 37232.25
$ws.Range("N141").Value = -47592.25

$ws = $wb.Worksheets.Item("LTW")
# Row 22
$ws.Range("H22").Value = 10500
$ws.Range("J22").Value = 10500
$ws.Range("L22").Value = 10500
$ws.Range("N22").Value = -11090

# Row 27
$ws.Range("H27").Value = 10500
$ws.Range("J27").Value = 10500
$ws.Range("L27").Value = 10500
$ws.Range("N27").Value = -10714

# Row 55
$ws.Range("H55").Value = 667100
$ws.Range("I55").Value = 1000175
$ws.Range("K55").Value = 1000175
$ws.Range("M55").Value = -1000002

# Row 62
$ws.Range("H62").Value = 40249
$ws.Range("J62").Value = 40249
$ws.Range("L62").Value = 40249
$ws.Range("N62").Value = -41497

# Row 65
$ws.Range("H65").Value = 40249
$ws.Range("J65").Value = 40249
$ws.Range("L65").Value = 120747
$ws.Range("N65").Value = -126987

# Row 95
$ws.Range("H95").Value = 28900
$ws.Range("J95").Value = 28900
$ws.Range("L95").Value = 28900
$ws.Range("N95").Value = -34392

$ws = $wb.Worksheets.Item("WVR")
# Row 135
$ws.Range("H135").Value = 56015.89
$ws.Range("J135").Value = 56015.89
$ws.Range("L135").Value = 56015.89
$ws.Range("N135").Value = -66155.89

# Row 141
$ws.Range("H141").Value = 45030.715
$ws.Range("J141").Value = 45030.715
$ws.Range("L141").Value = 45030.715
$ws.Range("N141").Value = -55390.715
